# Auto-generated edit script applying cell updates per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '244.73'

# Row 4
$ws.Range("B4").Value = 'HuobiToken'
$ws.Range("C4").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.030'
$ws.Range("E4").Value = '3HuobiTokenHT'

# Row 5
$ws.Range("B5").Value = 'Cronos'
$ws.Range("C5").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05615'
$ws.Range("E5").Value = '4CronosCRO'

# Row 6
$ws.Range("B6").Value = 'KuCoinToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.558'
$ws.Range("E6").Value = '5KuCoinTokenKCS'

# Row 7
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.010'
$ws.Range("E7").Value = '6GateTokenGT'

# Row 8
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8137'
$ws.Range("E8").Value = '7MXTokenMX'

# Row 9
$ws.Range("B9").Value = 'FTXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8374'
$ws.Range("E9").Value = '8FTXTokenFTT'

# Row 10
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1338'
$ws.Range("E10").Value = '9WazirXWRX'

# Row 11
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06957'
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'

# Row 12
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.02839'
$ws.Range("E12").Value = '11BitrueCoinBTR'

# Row 13
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09402'
$ws.Range("E13").Value = '12BitMartTokenBMX'

# Row 14
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001514'
$ws.Range("E14").Value = '13BitForexTokenBF'

# Row 15
$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0005986'
$ws.Range("E15").Value = '14OneONE'

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006190'

# Row 17
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.499'
$ws.Range("E17").Value = '16LEOLEO'

# Row 18
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.091'
$ws.Range("E18").Value = '17BTSETokenBTSE'

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3159'

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03255'

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1292'

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.741'

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04686'

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004530'

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.00009696'
$ws.Range("E27").Value = '26NitroExNTXBestin24h'

# Row 41
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1362'
$ws.Range("E41").Value = '40BKEXTokenBKK'

# Row 42
$ws.Range("B42").Value = 'KickToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.006288'
$ws.Range("E42").Value = '41KickTokenKICK'

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002734'

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008180'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005293'

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002016'

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001999'
